$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.335000000000001
$ws.Range("A4").Value = -20.872
$ws.Range("A6").Value = -22.074
$ws.Range("A7").Value = -21.486
$ws.Range("B7").Value = 6.355
$ws.Range("A8").Value = -21.559
$ws.Range("B11").Value = 6.709000000000001
$ws.Range("B12").Value = 5.739
$ws.Range("C12").Value = -12.13
$ws.Range("D12").Value = -7.631
$ws.Range("C13").Value = -13.332
$ws.Range("D13").Value = -8.327000000000002
$ws.Range("C14").Value = -12.241
$ws.Range("B15").Value = 5.087000000000001
$ws.Range("A16").Value = -21.795
$ws.Range("C16").Value = -13.17
$ws.Range("C19").Value = -12.635
$ws.Range("A20").Value = -22.116
$ws.Range("B20").Value = 5.413
$ws.Range("C20").Value = -12.365
$ws.Range("A21").Value = -21.709
$ws.Range("B21").Value = 6.692
$ws.Range("B22").Value = 8.124000000000001
$ws.Range("C22").Value = -13.29
$ws.Range("D22").Value = -8.063000000000001
$ws.Range("B23").Value = 8.094999999999999
$ws.Range("D25").Value = -8.339000000000002
$ws.Range("A28").Value = -21.749
$ws.Range("A29").Value = -21.648
$ws.Range("B29").Value = 6.181
$ws.Range("D29").Value = -7.171000000000001
$ws.Range("A30").Value = -21.465
$ws.Range("A32").Value = -21.486
$ws.Range("B34").Value = 7.696
$ws.Range("D34").Value = -8.026
$ws.Range("C36").Value = -12.89
$ws.Range("A40").Value = -20.745
$ws.Range("B42").Value = 7.392
$ws.Range("B43").Value = 4.513
$ws.Range("C43").Value = -13.13
$ws.Range("D43").Value = -8.452000000000002
$ws.Range("B44").Value = 6.152
$ws.Range("B45").Value = 5.166000000000001
$ws.Range("A46").Value = -20.864
$ws.Range("B46").Value = 7.3
$ws.Range("C46").Value = -13.484
$ws.Range("D48").Value = -8.058000000000002
$ws.Range("B50").Value = 5.1
$ws.Range("C50").Value = -13.419
$ws.Range("A51").Value = -21.134
$ws.Range("B51").Value = 6.7
$ws.Range("A52").Value = -21.783
$ws.Range("A57").Value = -21.157
$ws.Range("B57").Value = 7.263
$ws.Range("A59").Value = -22.134
$ws.Range("D60").Value = -8.360000000000001
$ws.Range("A62").Value = -22.087
$ws.Range("B65").Value = 5.038
$ws.Range("A66").Value = -21.64
$ws.Range("B66").Value = 6.422
$ws.Range("B67").Value = 5.8
$ws.Range("D68").Value = -7.025
$ws.Range("D70").Value = -7.389999999999999
$ws.Range("D71").Value = -7.784000000000001
$ws.Range("A73").Value = -20.767
$ws.Range("D73").Value = -8.396000000000001
$ws.Range("A74").Value = -21.244
$ws.Range("C76").Value = -12.136
$ws.Range("A77").Value = -21.176
$ws.Range("D78").Value = -8.172000000000001
$ws.Range("B79").Value = 5.499
$ws.Range("B84").Value = 5.962
$ws.Range("B87").Value = 5.048
$ws.Range("D87").Value = -8.294
$ws.Range("A92").Value = -21.472
$ws.Range("B92").Value = 5.499000000000001
$ws.Range("D92").Value = -6.25
$ws.Range("C95").Value = -11.862
$ws.Range("B97").Value = 7.041999999999999
$ws.Range("C97").Value = -13.281
$ws.Range("C99").Value = -12.123
$ws.Range("A100").Value = -21.481
$ws.Range("D101").Value = -8.390000000000001
